$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Reference / Payment reference values on row 2 from CR413367 to CR416232
$ws.Range("H2").Value = "CR416232"
$ws.Range("I2").Value = "CR416232"

# Update page setup to A4 / portrait (adds a <pageSetup> element to the sheet)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move/collapse the selection in the frozen (bottom-left) pane to cell A7
$ws.Range("A7").Select()
